# Adds a new parts row (SDRAM memory chip) to the "Spis czesci" sheet,
# mirroring the existing Altera MAX 10 row: a product-link cell, a name
# cell and a quantity of 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url  = "https://superelektronika.pl/pl/pamieci/247629409-k4d263238f-128m-ddr-sdram-pamiec-synchroniczna-1m-x-32bit-x-4-5905427014409.html"
$name = "K4D263238F"

# New data row right below the existing one (row 3 -> row 4).
$ws.Range("B4").Value = $url
$ws.Range("C4").Value = $name
$ws.Range("E4").Value = 1

# The link column wraps its (long) text, same as a manually entered URL.
$ws.Range("B4").WrapText = $true

# Row grows taller to fit the wrapped text.
$ws.Rows.Item(4).RowHeight = 28.8

# Keep the selection/active-cell in sync with the newly used range, same
# as Excel would after typing into the sheet (last edited cell D4).
$ws.Range("D4").Select()
